$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 11 (existing rows 11-15 shift down to 12-16).
# Excel copies the formatting of the row above (row 10) into the new row 11.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new "Baseline 2010 C82" data point.
$ws.Range("A11").Value = "CW3M"
$ws.Range("B11").Value = "Baseline 2010 C82"
$ws.Range("C11").Value = 2010
$ws.Range("D11").Value = 1090.199341
$ws.Range("E11").Value = 1990.4676509999999
$ws.Range("F11").Value = 1.255063
$ws.Range("G11").Value = 327.58108499999997
$ws.Range("H11").Value = 10.610913999999999
$ws.Range("I11").Value = 8.8404570000000007
$ws.Range("J11").Value = 814.49517800000001
$ws.Range("K11").Value = 93.229797000000005
$ws.Range("L11").Value = 1305.1243899999999
$ws.Range("M11").Value = 1201.781982
$ws.Range("N11").Value = 7126.6015630000002
$ws.Range("O11").Value = 29450.638672000001
$ws.Range("P11").Value = 3.3577499999999998
$ws.Range("Q11").Value = 0.00098200000000000002
$ws.Range("R11").Value = 2010

# Row 12 is the old placeholder row (previously row 11), which only had
# empty, style-only cells. After the insert it carries that same empty
# formatting across D:Q - clear it so only the C12 integer-style cell remains.
$ws.Range("D12:Q12").Clear()

# Move the active selection to B11, matching the edited workbook's saved view.
$ws.Range("B11").Select()
